# Auto-generated: update market-price columns (H-N) on several leve rows
# across the ALC/ARM/BSM/CRP/CUL/GSM sheets, per the scheduled price-refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
    # row 64: Forged from the Void | Void Glue
    $ws.Cells.Item(64, 8).Value = 6720
    $ws.Cells.Item(64, 9).Value = 6000
    $ws.Cells.Item(64, 10).Value = 6900
    $ws.Cells.Item(64, 11).Value = 6000
    $ws.Cells.Item(64, 12).Value = 6900
    $ws.Cells.Item(64, 13).Value = -5752
    $ws.Cells.Item(64, 14).Value = -7396
    # row 67: Dodging the Draft (L) | Void Glue
    $ws.Cells.Item(67, 8).Value = 6720
    $ws.Cells.Item(67, 9).Value = 6000
    $ws.Cells.Item(67, 10).Value = 6900
    $ws.Cells.Item(67, 11).Value = 6000
    $ws.Cells.Item(67, 12).Value = 6900
    $ws.Cells.Item(67, 13).Value = -5142
    $ws.Cells.Item(67, 14).Value = -8616
    # row 76: Warding Off Temptation | Enchanted Hardsilver Ink
    $ws.Cells.Item(76, 8).Value = 3507.6924
    $ws.Cells.Item(76, 9).Value = 3530
    $ws.Cells.Item(76, 10).Value = 3433.3333
    $ws.Cells.Item(76, 11).Value = 3530
    $ws.Cells.Item(76, 12).Value = 3433.3333
    $ws.Cells.Item(76, 13).Value = -3215
    $ws.Cells.Item(76, 14).Value = -4063.3333
    # row 79: The Garden of Arcane Delights (L) | Enchanted Hardsilver Ink
    $ws.Cells.Item(79, 8).Value = 3507.6924
    $ws.Cells.Item(79, 9).Value = 3530
    $ws.Cells.Item(79, 10).Value = 3433.3333
    $ws.Cells.Item(79, 11).Value = 3530
    $ws.Cells.Item(79, 12).Value = 3433.3333
    $ws.Cells.Item(79, 13).Value = -2438
    $ws.Cells.Item(79, 14).Value = -5617.3333
    # row 129: Practical Command | Commanding Craftsman's Draught
    $ws.Cells.Item(129, 8).Value = 811.61536
    $ws.Cells.Item(129, 9).Value = 536.2
    $ws.Cells.Item(129, 10).Value = 983.75
    $ws.Cells.Item(129, 11).Value = 1608.6
    $ws.Cells.Item(129, 12).Value = 2951.25
    $ws.Cells.Item(129, 13).Value = 3391.4
    $ws.Cells.Item(129, 14).Value = -12951.25

$ws = $wb.Worksheets.Item("ARM")
    # row 2: Ain't Got No Ingots | Bronze Ingot
    $ws.Cells.Item(2, 8).Value = 2273.2
    $ws.Cells.Item(2, 9).Value = 2318.25
    $ws.Cells.Item(2, 10).Value = 2231.6155
    $ws.Cells.Item(2, 11).Value = 2318.25
    $ws.Cells.Item(2, 12).Value = 2231.6155
    $ws.Cells.Item(2, 13).Value = -2205.25
    $ws.Cells.Item(2, 14).Value = -2457.6155
    # row 63: Rivets Run through It | Mythrite Rivets
    $ws.Cells.Item(63, 8).Value = 5069.1577
    $ws.Cells.Item(63, 9).Value = 5050.4165
    $ws.Cells.Item(63, 10).Value = 5101.2856
    $ws.Cells.Item(63, 11).Value = 5050.4165
    $ws.Cells.Item(63, 12).Value = 5101.2856
    $ws.Cells.Item(63, 13).Value = -4364.4165
    $ws.Cells.Item(63, 14).Value = -6473.2856
    # row 66: A Riveting Revival (L) | Mythrite Rivets
    $ws.Cells.Item(66, 8).Value = 5069.1577
    $ws.Cells.Item(66, 9).Value = 5050.4165
    $ws.Cells.Item(66, 10).Value = 5101.2856
    $ws.Cells.Item(66, 11).Value = 25252.0825
    $ws.Cells.Item(66, 12).Value = 25506.428
    $ws.Cells.Item(66, 13).Value = -21820.0825
    $ws.Cells.Item(66, 14).Value = -32370.428
    # row 116: No Scope | Titanbronze Ingot
    $ws.Cells.Item(116, 8).Value = 2273.2
    $ws.Cells.Item(116, 9).Value = 2318.25
    $ws.Cells.Item(116, 10).Value = 2231.6155
    $ws.Cells.Item(116, 11).Value = 2318.25
    $ws.Cells.Item(116, 12).Value = 2231.6155
    $ws.Cells.Item(116, 13).Value = -24.25
    $ws.Cells.Item(116, 14).Value = -6819.6155

$ws = $wb.Worksheets.Item("BSM")
    # row 3: Hells Bells | Bronze Ingot
    $ws.Cells.Item(3, 8).Value = 2273.2
    $ws.Cells.Item(3, 9).Value = 2318.25
    $ws.Cells.Item(3, 10).Value = 2231.6155
    $ws.Cells.Item(3, 11).Value = 2318.25
    $ws.Cells.Item(3, 12).Value = 2231.6155
    $ws.Cells.Item(3, 13).Value = -2204.25
    $ws.Cells.Item(3, 14).Value = -2459.6155
    # row 105: Ingot to Wing It | Molybdenum Ingot
    $ws.Cells.Item(105, 8).Value = 2762
    $ws.Cells.Item(105, 9).Value = 2874.2856
    $ws.Cells.Item(105, 10).Value = 2500
    $ws.Cells.Item(105, 11).Value = 2874.2856
    $ws.Cells.Item(105, 12).Value = 2500
    $ws.Cells.Item(105, 13).Value = -1127.2856
    $ws.Cells.Item(105, 14).Value = -5994

$ws = $wb.Worksheets.Item("CRP")
    # row 31: Wall Not Found | Walnut Lumber
    $ws.Cells.Item(31, 8).Value = 348510.12
    $ws.Cells.Item(31, 9).Value = 206162.4
    $ws.Cells.Item(31, 10).Value = 372234.72
    $ws.Cells.Item(31, 11).Value = 206162.4
    $ws.Cells.Item(31, 12).Value = 372234.72
    $ws.Cells.Item(31, 13).Value = -205867.4
    $ws.Cells.Item(31, 14).Value = -372824.72
    # row 34: Armoires of the Rich and Famous | Walnut Lumber
    $ws.Cells.Item(34, 8).Value = 348510.12
    $ws.Cells.Item(34, 9).Value = 206162.4
    $ws.Cells.Item(34, 10).Value = 372234.72
    $ws.Cells.Item(34, 11).Value = 206162.4
    $ws.Cells.Item(34, 12).Value = 372234.72
    $ws.Cells.Item(34, 13).Value = -205960.4
    $ws.Cells.Item(34, 14).Value = -372638.72
    # row 58: You Do the Heavy Lifting | Mahogany Lumber
    $ws.Cells.Item(58, 8).Value = 55557130
    $ws.Cells.Item(58, 9).Value = 55557130
    $ws.Cells.Item(58, 10).Value = 0
    $ws.Cells.Item(58, 11).Value = 55557130
    $ws.Cells.Item(58, 12).Value = 0
    $ws.Cells.Item(58, 13).Value = -55556927
    $ws.Cells.Item(58, 14).ClearContents()
    # row 62: Splinter in the Sewers | Cedar Lumber
    $ws.Cells.Item(62, 8).Value = 3502.4
    $ws.Cells.Item(62, 9).Value = 3666.6667
    $ws.Cells.Item(62, 10).Value = 3432
    $ws.Cells.Item(62, 11).Value = 3666.6667
    $ws.Cells.Item(62, 12).Value = 3432
    $ws.Cells.Item(62, 13).Value = -3042.6667
    $ws.Cells.Item(62, 14).Value = -4680
    # row 65: The Lumber of Their Discontent (L) | Cedar Lumber
    $ws.Cells.Item(65, 8).Value = 3502.4
    $ws.Cells.Item(65, 9).Value = 3666.6667
    $ws.Cells.Item(65, 10).Value = 3432
    $ws.Cells.Item(65, 11).Value = 18333.3335
    $ws.Cells.Item(65, 12).Value = 17160
    $ws.Cells.Item(65, 13).Value = -15213.3335
    $ws.Cells.Item(65, 14).Value = -23400
    # row 134: Wood You Be Quiet | Ceiba Lumber
    $ws.Cells.Item(134, 8).Value = 72918
    $ws.Cells.Item(134, 9).Value = 941.1111
    $ws.Cells.Item(134, 10).Value = 180883.33
    $ws.Cells.Item(134, 11).Value = 2823.3333
    $ws.Cells.Item(134, 12).Value = 542649.99
    $ws.Cells.Item(134, 13).Value = -288.3332999999998
    $ws.Cells.Item(134, 14).Value = -547719.99
    # row 135: The Wing's Wings | Ceiba Wings
    $ws.Cells.Item(135, 8).Value = 48514.285
    $ws.Cells.Item(135, 10).Value = 48514.285
    $ws.Cells.Item(135, 12).Value = 48514.285
    $ws.Cells.Item(135, 14).Value = -58654.285
    # row 136: Turali Quality | Dark Mahogany Lumber
    $ws.Cells.Item(136, 8).Value = 55557130
    $ws.Cells.Item(136, 9).Value = 55557130
    $ws.Cells.Item(136, 10).Value = 0
    $ws.Cells.Item(136, 11).Value = 166671390
    $ws.Cells.Item(136, 12).Value = 0
    $ws.Cells.Item(136, 13).Value = -166668840
    $ws.Cells.Item(136, 14).ClearContents()
    # row 137: Lament of the Lazylump | Dark Mahogany Fishing Rod
    $ws.Cells.Item(137, 8).Value = 40509.9
    $ws.Cells.Item(137, 10).Value = 40509.9
    $ws.Cells.Item(137, 12).Value = 40509.9
    $ws.Cells.Item(137, 14).Value = -50709.9
    # row 138: Bow Out | Acacia Longbow
    $ws.Cells.Item(138, 8).Value = 39960
    $ws.Cells.Item(138, 10).Value = 39960
    $ws.Cells.Item(138, 12).Value = 39960
    $ws.Cells.Item(138, 14).Value = -50240
    # row 139: Weaving a Path | Acacia Spinning Wheel
    $ws.Cells.Item(139, 8).Value = 53245.453
    $ws.Cells.Item(139, 10).Value = 53245.453
    $ws.Cells.Item(139, 12).Value = 53245.453
    $ws.Cells.Item(139, 14).Value = -63525.453
    # row 140: Spear Pressure | Claro Walnut Spear
    $ws.Cells.Item(140, 8).Value = 0
    $ws.Cells.Item(140, 10).Value = 0
    $ws.Cells.Item(140, 12).Value = 0
    $ws.Cells.Item(140, 14).ClearContents()
    # row 141: No Greater Treasure | Claro Walnut Necklace of Gathering
    $ws.Cells.Item(141, 8).Value = 52956.5
    $ws.Cells.Item(141, 10).Value = 52956.5
    $ws.Cells.Item(141, 12).Value = 52956.5
    $ws.Cells.Item(141, 14).Value = -63316.5

$ws = $wb.Worksheets.Item("CUL")
    # row 87: Soup That Eats Like a Knight | Clam Chowder
    $ws.Cells.Item(87, 8).Value = 26786.273
    $ws.Cells.Item(87, 9).Value = 11471.571
    $ws.Cells.Item(87, 10).Value = 33933.133
    $ws.Cells.Item(87, 11).Value = 34414.713
    $ws.Cells.Item(87, 12).Value = 101799.399
    $ws.Cells.Item(87, 13).Value = -33166.713
    $ws.Cells.Item(87, 14).Value = -104295.399
    # row 90: Like Ma Used to Make (L) | Clam Chowder
    $ws.Cells.Item(90, 8).Value = 26786.273
    $ws.Cells.Item(90, 9).Value = 11471.571
    $ws.Cells.Item(90, 10).Value = 33933.133
    $ws.Cells.Item(90, 11).Value = 103244.139
    $ws.Cells.Item(90, 12).Value = 305398.197
    $ws.Cells.Item(90, 13).Value = -97004.139
    $ws.Cells.Item(90, 14).Value = -317878.197
    # row 121: A Cookie for Your Troubles | Coffee Biscuit
    $ws.Cells.Item(121, 8).Value = 115256980
    $ws.Cells.Item(121, 9).Value = 1361.4286
    $ws.Cells.Item(121, 10).Value = 204900220
    $ws.Cells.Item(121, 11).Value = 4084.2858
    $ws.Cells.Item(121, 12).Value = 614700660
    $ws.Cells.Item(121, 13).Value = -2774.2858
    $ws.Cells.Item(121, 14).Value = -614703280

$ws = $wb.Worksheets.Item("GSM")
    # row 70: Sky Is the Limit | Mythrite Ingot
    $ws.Cells.Item(70, 8).Value = 37352.773
    $ws.Cells.Item(70, 9).Value = 52421.285
    $ws.Cells.Item(70, 11).Value = 52421.285
    $ws.Cells.Item(70, 13).Value = -52151.285
    # row 73: Hulls of Broken Dreams (L) | Mythrite Ingot
    $ws.Cells.Item(73, 8).Value = 37352.773
    $ws.Cells.Item(73, 9).Value = 52421.285
    $ws.Cells.Item(73, 11).Value = 52421.285
    $ws.Cells.Item(73, 13).Value = -51485.285
    # row 80: Needs More Prayerbell | Hardsilver Ingot
    $ws.Cells.Item(80, 8).Value = 3901.8333
    $ws.Cells.Item(80, 9).Value = 2205
    $ws.Cells.Item(80, 10).Value = 4056.0908
    $ws.Cells.Item(80, 11).Value = 2205
    $ws.Cells.Item(80, 12).Value = 4056.0908
    $ws.Cells.Item(80, 13).Value = -1207
    $ws.Cells.Item(80, 14).Value = -6052.0908
    # row 83: With a Noise That Reaches Heaven (L) | Hardsilver Ingot
    $ws.Cells.Item(83, 8).Value = 3901.8333
    $ws.Cells.Item(83, 9).Value = 2205
    $ws.Cells.Item(83, 10).Value = 4056.0908
    $ws.Cells.Item(83, 11).Value = 11025
    $ws.Cells.Item(83, 12).Value = 20280.454
    $ws.Cells.Item(83, 13).Value = -6033
    $ws.Cells.Item(83, 14).Value = -30264.454

